$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row of data: Location / e4W4rmGe9QzuGm2Dy4NBqVc0KDe6yGld6HW95UuN-Qd03
# (set B7 first so the shared-strings table order matches: UUID string before "Location")
$ws.Range("B7").Value = "e4W4rmGe9QzuGm2Dy4NBqVc0KDe6yGld6HW95UuN-Qd03"
$ws.Range("A7").Value = "Location"

# Update the active cell selection, matching the author's final cursor position
$ws.Range("B13").Select()
